# Fruta / hortaliza, semanal
# Insert two new weekly records (Early Burlat / Santina, fair date 2021-11-22)
# ahead of the existing row 78 ("Bing" / "Primera", 2020-12-01), pushing the
# rest of the Cereza - Macroferia Regional de Talca rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 78; everything from the old row 78 onward
# (through row 88) shifts down to rows 80-90.
$ws.Rows.Item(78).Insert()
$ws.Rows.Item(78).Insert()

# New row 78: Early Burlat / Primera
$ws.Cells.Item(78, 1).Value = 5
$ws.Cells.Item(78, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(78, 3).Value = "Maule"
$ws.Cells.Item(78, 4).Value = 44522
$ws.Cells.Item(78, 5).Value = 7
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100103
$ws.Cells.Item(78, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(78, 9).Value = 100103001
$ws.Cells.Item(78, 10).Value = "Cereza"
$ws.Cells.Item(78, 11).Value = "Early Burlat"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 50
$ws.Cells.Item(78, 14).Value = 30000
$ws.Cells.Item(78, 15).Value = 30000
$ws.Cells.Item(78, 16).Value = 30000
$ws.Cells.Item(78, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(78, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(78, 19).Value = 2000
$ws.Cells.Item(78, 20).Value = 15

# New row 79: Santina / Primera
$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = 44522
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100103
$ws.Cells.Item(79, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(79, 9).Value = 100103001
$ws.Cells.Item(79, 10).Value = "Cereza"
$ws.Cells.Item(79, 11).Value = "Santina"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 60
$ws.Cells.Item(79, 14).Value = 20000
$ws.Cells.Item(79, 15).Value = 20000
$ws.Cells.Item(79, 16).Value = 20000
$ws.Cells.Item(79, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(79, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(79, 19).Value = 2000
$ws.Cells.Item(79, 20).Value = 10

# Make sure the date column keeps the same number format as its neighbours
# (the insert already copies formatting from the row below, but set it
# explicitly so both new rows match column D's date style).
$dateFormat = $ws.Range("D80").NumberFormat
$ws.Range("D78").NumberFormat = $dateFormat
$ws.Range("D79").NumberFormat = $dateFormat
